$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell values for the two new data rows.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "MCH218-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24E | GRAP COUNT NUMER: NONE"

$ws.Range("A3").Value = "MCH218-2"
$ws.Range("C3").Value = "DELIGATES"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 24E | GRAP COUNT NUMER: NONE"

# ---------------------------------------------------------------------
# 2. Formatting: build the 10pt Calibri / automatic-text-color style on a
#    scratch cell once, then stamp it onto every data cell via
#    Copy + PasteSpecial(xlPasteFormats) so we don't fork a brand-new
#    style per target cell.
# ---------------------------------------------------------------------
$seed = $ws.Range("Z100")
$seed.Font.Name = "Calibri"
$seed.Font.Size = 10
$seed.Font.ThemeColor = 1
$seed.Copy()

foreach ($addr in @("A2", "C2:D2", "E2:H2", "A3", "C3:D3", "E3:H3")) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}

# F2/F3 ("extentAndMedium") use a second, alignment-flagged variant of the
# same font - build it once and stamp it the same way.
$seed2 = $ws.Range("Z101")
$seed2.Font.Name = "Calibri"
$seed2.Font.Size = 10
$seed2.Font.ThemeColor = 1
$seed2.WrapText = $false
$seed2.Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F3").PasteSpecial(-4122)

$ws.Range("Z100:Z101").Clear()
$excel.CutCopyMode = $false
